# Textbox response formatting fix
# Renames the task-order sheets and refreshes the generated stimulus
# file names/timestamps that are written into column B of each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "GNG_TO-1651168757207226"
$ws.Range("B2").Value = "go_stims-1651168757170438.csv"
$ws.Range("B3").Value = "GNG_stims-16511687571884568.csv"
$ws.Range("B4").Value = "go_stims-16511687571884568.csv"
$ws.Range("B5").Value = "GNG_stims-16511687572052257.csv"

# --- Sheet 2: NB_TO ------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-16511687600817351"
$ws.Range("B2").Value = "TB-16511687600596895.csv"
$ws.Range("B3").Value = "TB-16511687597424212.csv"
$ws.Range("B4").Value = "OB-16511687586097894.csv"
$ws.Range("B5").Value = "ZB-match_1-16511687573053324.csv"
$ws.Range("B6").Value = "ZB-match_9-16511687574727373.csv"
$ws.Range("B7").Value = "OB-16511687579176607.csv"
$ws.Range("B8").Value = "OB-16511687577229903.csv"
$ws.Range("B9").Value = "TB-165116875950569.csv"
$ws.Range("B10").Value = "ZB-match_4-1651168757388458.csv"

# --- Sheet 3: RS_TO (name only, cell contents unchanged) -----------------
$ws = $wb.Worksheets.Item(3)
$ws.Name = "RS_TO-16511687600837328"

# --- Sheet 4: TOL_TO -------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Name = "TOL_TO-16511687601308706"
$ws.Range("B2").Value = "MM_stims-16511687600976956.csv"
$ws.Range("B3").Value = "ZM_stims-16511687600848298.csv"
$ws.Range("B4").Value = "MM_stims-1651168760114555.csv"
$ws.Range("B5").Value = "ZM_stims-16511687600976956.csv"
$ws.Range("B6").Value = "MM_stims-1651168760129908.csv"
$ws.Range("B7").Value = "ZM_stims-16511687601155572.csv"

# --- Sheet 5: vSAT_TO -------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Name = "vSAT_TO-16511687602141945"
$ws.Range("B2").Value = "vSAT_stims-1651168760179355.csv"
$ws.Range("B3").Value = "SAT_stims-16511687601633203.csv"
$ws.Range("B4").Value = "SAT_stims-16511687601358035.csv"
$ws.Range("B5").Value = "vSAT_stims-16511687601985297.csv"
